$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing values for row 23 (C, D, E columns)
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 5

# Update the active selection to G23
$ws.Range("G23").Select()
